# Auto-generated edit script applying the diff to Lich_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1045.5264
$ws.Range("J19").Value = 1398.2727
$ws.Range("L19").Value = 1398.2727
$ws.Range("N19").Value = -1748.2727
$ws.Range("H39").Value = 805.46155
$ws.Range("I39").Value = 677.4
$ws.Range("K39").Value = 2032.2
$ws.Range("M39").Value = -1736.2
$ws.Range("H62").Value = 57697090
$ws.Range("J62").Value = 83338430
$ws.Range("L62").Value = 83338430
$ws.Range("N62").Value = -83339678
$ws.Range("H65").Value = 57697090
$ws.Range("J65").Value = 83338430
$ws.Range("L65").Value = 416692150
$ws.Range("N65").Value = -416698390
$ws.Range("H74").Value = 3595.6
$ws.Range("I74").Value = 4199.3335
$ws.Range("K74").Value = 4199.3335
$ws.Range("M74").Value = -3263.3335
$ws.Range("H77").Value = 3595.6
$ws.Range("I77").Value = 4199.3335
$ws.Range("K77").Value = 20996.6675
$ws.Range("M77").Value = -16316.6675
$ws.Range("H86").Value = 56940.684
$ws.Range("I86").Value = 103868.2
$ws.Range("J86").Value = 4799
$ws.Range("K86").Value = 103868.2
$ws.Range("L86").Value = 4799
$ws.Range("M86").Value = -102745.2
$ws.Range("N86").Value = -7045
$ws.Range("H89").Value = 56940.684
$ws.Range("I89").Value = 103868.2
$ws.Range("J89").Value = 4799
$ws.Range("K89").Value = 519341
$ws.Range("L89").Value = 23995
$ws.Range("M89").Value = -513725
$ws.Range("N89").Value = -35227
$ws.Range("H112").Value = 4980880.5
$ws.Range("J112").Value = 5363648.5
$ws.Range("L112").Value = 16090945.5
$ws.Range("N112").Value = -16093161.5
$ws.Range("H132").Value = 1591.3948
$ws.Range("I132").Value = 1431.7028
$ws.Range("K132").Value = 4295.1084
$ws.Range("M132").Value = -1765.1084
$ws.Range("H135").Value = 1377.6666
$ws.Range("I135").Value = 1015.8823
$ws.Range("K135").Value = 9142.940699999999
$ws.Range("M135").Value = -6607.940699999999
$ws.Range("H137").Value = 2853.8
$ws.Range("I137").Value = 2893.718
$ws.Range("J137").Value = 2803.5806
$ws.Range("K137").Value = 8681.153999999999
$ws.Range("L137").Value = 8410.7418
$ws.Range("M137").Value = -6131.153999999999
$ws.Range("N137").Value = -13510.7418
$ws.Range("H138").Value = 1495.3334
$ws.Range("I138").Value = 726.8293
$ws.Range("J138").Value = 2038.5862
$ws.Range("K138").Value = 2180.4879
$ws.Range("L138").Value = 6115.7586
$ws.Range("M138").Value = 2959.5121
$ws.Range("N138").Value = -16395.7586
$ws.Range("H141").Value = 1966.3636
$ws.Range("I141").Value = 1190.5186
$ws.Range("K141").Value = 3571.5558
$ws.Range("M141").Value = 1608.4442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2698.1606
$ws.Range("I61").Value = 1683.6097
$ws.Range("K61").Value = 1683.6097
$ws.Range("M61").Value = -1471.6097
$ws.Range("H132").Value = 2550.24
$ws.Range("I132").Value = 2712.932
$ws.Range("J132").Value = 1357.1666
$ws.Range("K132").Value = 8138.795999999999
$ws.Range("L132").Value = 4071.4998
$ws.Range("M132").Value = -5608.795999999999
$ws.Range("N132").Value = -9131.4998
$ws.Range("H136").Value = 2698.1606
$ws.Range("I136").Value = 1683.6097
$ws.Range("K136").Value = 5050.8291
$ws.Range("M136").Value = -2500.8291

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 979.3333
$ws.Range("I107").Value = 1298.5714
$ws.Range("K107").Value = 1298.5714
$ws.Range("M107").Value = 621.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1880.3077
$ws.Range("I16").Value = 1557.2858
$ws.Range("K16").Value = 1557.2858
$ws.Range("M16").Value = -1270.2858
$ws.Range("H31").Value = 71484856
$ws.Range("I31").Value = 2993.8
$ws.Range("J31").Value = 111197000
$ws.Range("K31").Value = 2993.8
$ws.Range("L31").Value = 111197000
$ws.Range("M31").Value = -2698.8
$ws.Range("N31").Value = -111197590
$ws.Range("H34").Value = 71484856
$ws.Range("I34").Value = 2993.8
$ws.Range("J34").Value = 111197000
$ws.Range("K34").Value = 2993.8
$ws.Range("L34").Value = 111197000
$ws.Range("M34").Value = -2791.8
$ws.Range("N34").Value = -111197404
$ws.Range("H48").Value = 14999
$ws.Range("I48").Value = 14999
$ws.Range("K48").Value = 14999
$ws.Range("M48").Value = -14523
$ws.Range("H113").Value = 1880.3077
$ws.Range("I113").Value = 1557.2858
$ws.Range("K113").Value = 1557.2858
$ws.Range("M113").Value = 612.7141999999999
$ws.Range("H132").Value = 1081.409
$ws.Range("I132").Value = 1081.409
$ws.Range("K132").Value = 3244.227
$ws.Range("M132").Value = -714.2270000000003
$ws.Range("H134").Value = 1069.7273
$ws.Range("I134").Value = 1132.2439
$ws.Range("K134").Value = 3396.7317
$ws.Range("M134").Value = -861.7316999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3692888
$ws.Range("J4").Value = 70796.92999999999
$ws.Range("L4").Value = 212390.79
$ws.Range("N4").Value = -212614.79
$ws.Range("H26").Value = 1836.35
$ws.Range("I26").Value = 195.18182
$ws.Range("J26").Value = 3842.2222
$ws.Range("K26").Value = 585.5454599999999
$ws.Range("L26").Value = 11526.6666
$ws.Range("M26").Value = -297.5454599999999
$ws.Range("N26").Value = -12102.6666
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
$ws.Range("H122").Value = 976.9
$ws.Range("J122").Value = 1021.875
$ws.Range("L122").Value = 9196.875
$ws.Range("N122").Value = -14096.875
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3622.125
$ws.Range("I97").Value = 3241.476
$ws.Range("K97").Value = 3241.476
$ws.Range("M97").Value = -2745.476

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 38462200
$ws.Range("I16").Value = 38462200
$ws.Range("K16").Value = 38462200
$ws.Range("M16").Value = -38462030
$ws.Range("H132").Value = 3831.4927
$ws.Range("I132").Value = 2958.14
$ws.Range("K132").Value = 8874.42
$ws.Range("M132").Value = -6344.42
